# Applies the diff: add jc="both" (justify) to 13 paragraphs in the
# Iterator/Mediator sections, and append a new "Memento (1)" section
# after the final "Colleague" paragraph.

$d = $word.ActiveDocument

# --- Part 1: add justify alignment to the 13 existing paragraphs -----
# (Iterator details: paragraphs 106-112, Mediator details: 114-119;
#  the two section headings 105 and 113 are left untouched.)
$justifyIdx = @(106,107,108,109,110,111,112,114,115,116,117,118,119)
foreach ($i in $justifyIdx) {
    $d.Paragraphs($i).Alignment = 3   # wdAlignParagraphJustify
}

# --- Part 2: append the new "Memento (1)" section ---------------------

# Heading2 "Memento (1)"
$tail = $d.Paragraphs($d.Paragraphs.Count).Range
$tail.Collapse(0)
$tail.InsertParagraphAfter()
$p = $d.Paragraphs($d.Paragraphs.Count)
$p.Range.ListFormat.RemoveNumbers()
$p.Range.Text = "Memento (1)"
$p.Style = "Heading 2"

# Paragraph: "Memento patern ... " + italic "memento" + "."
$tail = $d.Paragraphs($d.Paragraphs.Count).Range
$tail.Collapse(0)
$tail.InsertParagraphAfter()
$p = $d.Paragraphs($d.Paragraphs.Count)
$p.Range.Text = "Memento patern, resava problem narusene enkapsulacije kod pokusaja da sacuvamo prethodno stanje nekog objekta. Memento patern ustvari delegira pravljenje snepshota stanja objektu cije se stanje i zeli zapamtiti u memento – originatoru. Tako da umesto da se pokusava zapamtiti snaphsot objekta sa spoljasnje strane, od strane nekog drugog objekta i time se naruse svi skriveni fildovi i metode, sam originator pravi svoj snepshot i smesta ga u specijalan objekat "
$p.Style = "Normal"
$p.Alignment = 3
$endRng = $p.Range
$endRng.Collapse(0)
$endRng.InsertAfter("memento")
$endRng.Font.Italic = $true
$endRng.Font.ItalicBi = $true
$endRng2 = $p.Range
$endRng2.Collapse(0)
$endRng2.InsertAfter(".")
$endRng2.Font.Italic = $false
$endRng2.Font.ItalicBi = $false

# Paragraph: "Sadrzaj memento objekta..."
$tail = $d.Paragraphs($d.Paragraphs.Count).Range
$tail.Collapse(0)
$tail.InsertParagraphAfter()
$p = $d.Paragraphs($d.Paragraphs.Count)
$p.Range.Text = "Sadrzaj memento objekta nije pristupacan ni jednom drugom objektu osim onom koji ga je kreirao. Ostali objekti mogu, uz pomoc limitiranog interfejsa, samo da procitaju metadata mementa (creation time, name of operation...) ali ne detaljne vrednosti samog stanja objekta sadrzanog u snepshotu."
$p.Style = "Normal"
$p.Alignment = 3

# Paragraph: "Tako kreirani mementoi..."
$tail = $d.Paragraphs($d.Paragraphs.Count).Range
$tail.Collapse(0)
$tail.InsertParagraphAfter()
$p = $d.Paragraphs($d.Paragraphs.Count)
$p.Range.Text = "Tako kreirani mementoi se cuvaju u posebnim objektima nazvanim caretakers. Caretaker radi sa mementom samo preko posebnog limitiranog interfejsa, ne moze da pristupi detaljima snepshota niti da ih menja. U isto vreme, samo originator ima pristup poljima u mementu i to mu omogucava da se vrati u svoje predjasnje stanje."
$p.Style = "Normal"
$p.Alignment = 3

# Paragraph: "Klase koje ucestvuju:"
$tail = $d.Paragraphs($d.Paragraphs.Count).Range
$tail.Collapse(0)
$tail.InsertParagraphAfter()
$p = $d.Paragraphs($d.Paragraphs.Count)
$p.Range.Text = "Klase koje ucestvuju:"
$p.Style = "Normal"
$p.Alignment = 3

# List paragraph: "Memento - cuva interno stanje..."
$tail = $d.Paragraphs($d.Paragraphs.Count).Range
$tail.Collapse(0)
$tail.InsertParagraphAfter()
$p = $d.Paragraphs($d.Paragraphs.Count)
$p.Range.Text = "Memento – cuva interno stanje Originator objekta. Moze da cuva samo one podatke koje originator odluci da su bitni za cuvanje. Stiti podatke od pristupa drugih objekata koji nisu originatori. Memento efektivno ima dva interfejsa, onaj limitirani koji vidi caretaker i onaj potpuni koji vidi Originator kako bi imao pun pristup."
$p.Style = "List Paragraph"
$p.Alignment = 3
$p.Range.ListFormat.ApplyListTemplateWithLevel($d.Paragraphs(119).Range.ListFormat.ListTemplate)

# List paragraph: "Originator - kreira memento..."
$tail = $d.Paragraphs($d.Paragraphs.Count).Range
$tail.Collapse(0)
$tail.InsertParagraphAfter()
$p = $d.Paragraphs($d.Paragraphs.Count)
$p.Range.Text = "Originator – kreira memento sa sadrzanim svojim snepshotom, koristi vec kreiran memento da vrati svoje predjasnje stanje"
$p.Style = "List Paragraph"
$p.Alignment = 3

# List paragraph: "Caretaker" + " - odgovoran za cuvanje mementa..."
$tail = $d.Paragraphs($d.Paragraphs.Count).Range
$tail.Collapse(0)
$tail.InsertParagraphAfter()
$p = $d.Paragraphs($d.Paragraphs.Count)
$p.Range.Text = "Caretaker"
$p.Style = "List Paragraph"
$p.Alignment = 3
$endRng = $p.Range
$endRng.Collapse(0)
$endRng.InsertAfter(" – odgovoran za cuvanje mementa, nikad ne pristupa detaljima mementa.")

Write-Host "Done. Paragraph count: $($d.Paragraphs.Count)"
